$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header H1: pipe_length -> length_pipe
$ws.Range("H1").Value = "length_pipe"

# Update numeric values in row 2
$ws.Range("B2").Value = [double]"6.129859035316447e-07"
$ws.Range("G2").Value = 3.468721592776412
$ws.Range("K2").Value = 0.0196
